$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flattened new values for range G2:T26 (25 rows x 14 cols), row-major order
# Columns: G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$flat = @(44.04223000000001,132.12669,0.1792438957545786,0.1868246872369915,3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,295.6666326741901,2660.99969406771,0.05559858971965417,0.06037505722158815,44.04223000000001,132.12669,0.1792438957545786,0.1868246872369915,3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,14.72540215455334,132.52861939098,0.002769036145347189,0.003006923674988178,44.04223000000001,132.12669,0.1792438957545786,0.1868246872369915,3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,293.3998524998067,2640.598672498261,0.0551723333654613,0.05991218124028672,44.04223000000001,132.12669,0.1792438957545786,0.1868246872369915,2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,114.85852437714,689.15114626284,0.02159855481455877,0.01563605587129515,44.04223000000001,132.12669,0.1792438957545786,0.1868246872369915,3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,234.5471307568934,2110.924176812041,0.04410538170955722,0.04789446922883331,49.58946566666666,148.768397,0.2018201397722426,0.2103557520458098,3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,332.9066291551469,2996.159662396322,0.06260137953999778,0.06797945579079397,49.58946566666666,148.768397,0.2018201397722426,0.2103557520458098,3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,16.58010560707489,149.220950463674,0.003117803591222639,0.003385653686089768,49.58946566666666,148.768397,0.2018201397722426,0.2103557520458098,3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,330.3543420063931,2973.189078057538,0.06212143506758014,0.06745827935211973,49.58946566666666,148.768397,0.2018201397722426,0.2103557520458098,2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,129.325260122482,775.9515607348919,0.02431895007192369,0.01760545857445621,49.58946566666666,148.768397,0.2018201397722426,0.2103557520458098,3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,264.0889638849835,2376.800674964852,0.04966057150151832,0.05392690464235012,51.56497066666667,154.694912,0.2098600871655206,0.2187357342528896,3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,346.168694022912,3115.518246206208,0.06509524263421726,0.07068756633416415,51.56497066666667,154.694912,0.2098600871655206,0.2187357342528896,3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,17.24061043581156,155.165493922304,0.003242008127421512,0.003520528617594316,51.56497066666667,154.694912,0.2098600871655206,0.2187357342528896,3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,343.5147309243165,3091.632578318849,0.0645961785223311,0.07014562769031908,51.56497066666667,154.694912,0.2098600871655206,0.2187357342528896,2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,134.477215171072,806.863291026432,0.02528774872333019,0.01830681058487946,51.56497066666667,154.694912,0.2098600871655206,0.2187357342528896,3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,274.6095262985103,2471.485736686592,0.05163890915822052,0.05607520102593257,29.9106925,59.821385,0.1217310987299521,0.08458632803643724,3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,200.7980461577025,1204.788276946215,0.03775903991551374,0.0273352760327959,29.9106925,59.821385,0.1217310987299521,0.08458632803643724,3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,10.00056027552833,60.00336165317,0.001880553928919246,0.001361408045770938,29.9106925,59.821385,0.1217310987299521,0.08458632803643724,3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,199.2585926658817,1195.55155599529,0.03746955360348018,0.0271257053375436,29.9106925,59.821385,0.1217310987299521,0.08458632803643724,2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,78.00463336321499,312.01853345286,0.01466837014162684,0.007079345726122841,29.9106925,59.821385,0.1217310987299521,0.08458632803643724,3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,159.2895524324433,955.73731459466,0.02995358114041206,0.02168459289420396,70.603826,211.811478,0.2873447785777061,0.2994974984278718,3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,473.9813467059779,4265.832120353802,0.08912975465619818,0.09678688011059049,70.603826,211.811478,0.2873447785777061,0.2994974984278718,3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,23.60620094623066,212.455808516076,0.004439024685939009,0.004820380710607655,70.603826,211.811478,0.2873447785777061,0.2994974984278718,3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,470.3474854548013,4233.127369093212,0.08844642573614059,0.09604484649323312,70.603826,211.811478,0.2873447785777061,0.2994974984278718,2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,184.128988694268,1104.773932165608,0.0346245093851644,0.02506606427656368,70.603826,211.811478,0.2873447785777061,0.2994974984278718,3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,376.0010519167386,3384.009467250648,0.07070506411426397,0.07677932683687679)

$nRows = 25
$nCols = 14
$data = New-Object "object[,]" $nRows,$nCols
$idx = 0
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $data[$r,$c] = $flat[$idx]
        $idx++
    }
}

$ws.Range("G2:T26").Value2 = $data
